$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Table row "6" (row 7, col 1) gains four new trailing paragraphs:
#    three whitespace-only paragraphs and one "...Done" paragraph,
#    mirroring the pattern already used by the other rows in the plan.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 1)
$firstPara = $cell.Range.Paragraphs.Item(1)

$sp1 = "".PadLeft(53)
$sp2 = "".PadLeft(83)
$sp3 = "".PadLeft(97)
$sp4 = "".PadLeft(74) + "Done"

$newParas = "`r" + $sp1 + "`r" + $sp2 + "`r" + $sp3 + "`r" + $sp4
$firstPara.Range.InsertAfter($newParas)

# ------------------------------------------------------------------
# 2) Collapse the "(Bkash ,wallet)" run-split (with proofErr markers
#    from a spelling/grammar check) back into a single plain run.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(Bkash ,wallet)", $true, $false, $false, $false, $false, $true, 1, $false, "(Bkash ,wallet)", 2) | Out-Null
